# edit.ps1 -- applies the authored changes to Presentacion2tri.pptx
#
# Summary of the edit (see commit message / xml diff):
#  - Slide 1  : credits list, last line "Graciela Arias" gets a co-author appended.
#  - Slide 2  : slide is hidden from the slide show.
#  - Slide 4  : slide is hidden from the slide show.
#  - Slide 6  : slide is hidden from the slide show.
#  - Slide 10 : slide is shown again (was hidden); title textbox is repositioned
#               and its text is shortened.
#  - Slide 7  : title textbox ("Consultas DDL.") repositioned, reworded and
#               its run language switched to es-CO.
#  - Slide 8  : title textbox ("Consultas DML.") repositioned, reworded and
#               its run language switched to es-CO.

$p = $ppt.ActivePresentation

# Helper: point <-> EMU conversion (PowerPoint COM positions/sizes are in
# points; the underlying OOXML stores English Metric Units, 12700 EMU/pt).
# A tiny epsilon is added before handing the value to the COM layer so that
# the point -> EMU round trip lands exactly on the target EMU instead of
# being floored to one EMU below it.
function EmuToPt([double]$emu) {
    return ($emu / 12700.0) + 0.000015
}

# ---------------------------------------------------------------------------
# Slide 1 - "Graciela Arias" -> "Graciela Arias - Maria Pilar Bonilla"
# ---------------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$creditsShape = $s1.Shapes.Item(2)
$creditsRange = $creditsShape.TextFrame.TextRange
$oldName = "Graciela Arias"
$newName = "Graciela Arias – María Pilar Bonilla"
$fullText = $creditsRange.Text
$pos = $fullText.IndexOf($oldName)
if ($pos -ge 0) {
    $hit = $creditsRange.Characters($pos + 1, $oldName.Length)
    $hit.Text = $newName
}

# ---------------------------------------------------------------------------
# Slides 2, 4, 6 - hide from slide show
# ---------------------------------------------------------------------------
$p.Slides.Item(2).SlideShowTransition.Hidden = $true
$p.Slides.Item(4).SlideShowTransition.Hidden = $true
$p.Slides.Item(6).SlideShowTransition.Hidden = $true

# ---------------------------------------------------------------------------
# Slide 10 - un-hide, move/resize title, retitle, drop the autofit shrink
# ---------------------------------------------------------------------------
$s10 = $p.Slides.Item(10)
$s10.SlideShowTransition.Hidden = $false

$title10 = $s10.Shapes.Item(2)
$title10.Left = EmuToPt 2788364
$title10.Top  = EmuToPt 2316479
$title10.TextFrame.AutoSize = 2
$title10.TextFrame.TextRange.Text = "Inventario"

# ---------------------------------------------------------------------------
# Slide 7 - "Consultas DDL." -> "Lenguaje de definición de datos (DDL)"
# ---------------------------------------------------------------------------
$title7 = $p.Slides.Item(7).Shapes.Item(2)
$title7.Left = EmuToPt 3525810
$title7.Top  = EmuToPt 3207635
$title7.TextFrame.TextRange.Text = "Lenguaje de definición de datos (DDL)"
$title7.TextFrame.TextRange.LanguageID = "es-CO"

# ---------------------------------------------------------------------------
# Slide 8 - "Consultas DML." -> "Lenguaje de manipulación de datos (DML)"
# ---------------------------------------------------------------------------
$title8 = $p.Slides.Item(8).Shapes.Item(2)
$title8.Left = EmuToPt 3525810
$title8.Top  = EmuToPt 2997550
$title8.TextFrame.TextRange.Text = "Lenguaje de manipulación de datos (DML)"
$title8.TextFrame.TextRange.LanguageID = "es-CO"
